$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.710.49'
$ws.Range("E2").Value = '  +1.42%  '
$ws.Range("D3").Value = '3.163.83'
$ws.Range("E3").Value = '  +1.28%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '615.86'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +2.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.08'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -2.33%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '3.159.83'
$ws.Range("E8").Value = '  +1.26%  '
$ws.Range("E9").Value = '  -0.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.152'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -0.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.52'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  -1.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.474'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  -0.85%  '
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.83'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -2.69%  '
$ws.Range("D15").Value = '3.681.28'
$ws.Range("E15").Value = '  +1.88%  '
$ws.Range("E16").Value = '  +3.02%  '
$ws.Range("D17").Value = '64.680.38'
$ws.Range("E17").Value = '  +1.03%  '
$ws.Range("D18").Value = '3.161.87'
$ws.Range("E18").Value = '  +0.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.92'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -1.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '478.66'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -0.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.70'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("E22").Value = '  +1.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.97'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +2.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.77'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.82'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +0.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.83'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -3.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.58'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("E29").Value = '  -6.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.89'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -1.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.09'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -6.95%  '
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.69'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -0.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.63'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -0.18%  '
$ws.Range("E35").Value = '  +2.97%  '
$ws.Range("D36").Value = '0.0₃0789'
$ws.Range("E36").Value = '  +5.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.02'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -1.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.21'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -1.33%  '
$ws.Range("E39").Value = '  -2.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '465.02'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +3.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0400'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +0.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.120'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -3.61%  '
$ws.Range("E43").Value = '  -1.04%  '
$ws.Range("D44").Value = '2.848.90'
$ws.Range("E44").Value = '  -1.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.33'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -0.81%  '
$ws.Range("E46").Value = '  -1.41%  '
$ws.Range("E47").Value = '  +5.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.69'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -0.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.999'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("E50").Value = '  -1.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.71'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +1.04%  '
